# Generate Report for Handoff
#
# Two new source files were processed by the localization pipeline since the
# last report was generated:
#   5c4cc5a0-b7f6-4851-8ce5-f381df8f46f4.md  (Ready for handoff)
#   c97306d9-c77a-49c6-abdf-eea21385d93f.md  (Ready for handoff)
#
# The first sorts alphabetically/chronologically between the existing
# "d4a6720e..." and "ba776950..." entries, so it is inserted as a new row
# ahead of "ba776950...", which shifts down by one row. The second is newer
# still and is appended as the last row.
#
# This script rewrites the three worksheets (Overview / zh-cn / de-de) so the
# final row order/content reflects that, and rebuilds every hyperlink on each
# sheet (COM hyperlink mutation/removal of individual links is unreliable in
# this host, but bulk delete + re-add is not, so that's the approach used).

$wb = $excel.ActiveWorkbook

function Add-Link($ws, $cellRef, $url, $text) {
    $ws.Hyperlinks.Add($ws.Range($cellRef), $url, [System.Type]::Missing, [System.Type]::Missing, $text) | Out-Null
}

# ---------------------------------------------------------------------
# Sheet 1: "Overview"
# Columns: A File Name | B zh-cn | C de-de | D Latest Handoff Date
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overview")

$ov = @(
    @{ Row=2; Uuid="8e57db03-c7a5-42e7-aa11-fa025c83e4e3.md"; ZhCn="Handed back: in sync with en-US"; DeDe="Handed back: in sync with en-US"; Date="2016-25-12 02:25:41" },
    @{ Row=3; Uuid="65659143-8ae7-4026-91aa-8a2aa6f21603.md"; ZhCn="In Translation";                  DeDe="In Translation";                  Date="2016-26-12 02:26:26" },
    @{ Row=4; Uuid="d4a6720e-54cf-49b4-a943-9d2147c24d82.md"; ZhCn="In Translation";                  DeDe="In Translation";                  Date="2016-26-12 02:26:26" },
    @{ Row=5; Uuid="5c4cc5a0-b7f6-4851-8ce5-f381df8f46f4.md"; ZhCn="Ready for handoff";                DeDe="Ready for handoff";                Date="2016-26-12 02:26:59" },
    @{ Row=6; Uuid="ba776950-b98c-4ead-8535-7f56a0869957.md"; ZhCn="Ready for handoff";                DeDe="Ready for handoff";                Date="2016-25-12 02:25:19" },
    @{ Row=7; Uuid="c97306d9-c77a-49c6-abdf-eea21385d93f.md"; ZhCn="Ready for handoff";                DeDe="Ready for handoff";                Date="2016-26-12 02:26:59" }
)

foreach ($r in $ov) {
    $row = $r.Row
    $ws1.Range("A$row").Value2 = $r.Uuid
    $ws1.Range("B$row").Value2 = $r.ZhCn
    $ws1.Range("C$row").Value2 = $r.DeDe
    $ws1.Range("D$row").Value2 = $r.Date
}

# Rebuild hyperlinks (column A only) from scratch in final row order.
$ws1.Hyperlinks.Delete() | Out-Null

Add-Link $ws1 "A2" "https://github.com/OpenLocalizationTest/oltest/blob/07144325e6973bba9c1363117d82640c93e928c3/e2e/8e57db03-c7a5-42e7-aa11-fa025c83e4e3.md" "8e57db03-c7a5-42e7-aa11-fa025c83e4e3.md"
Add-Link $ws1 "A3" "https://github.com/OpenLocalizationTest/oltest/blob/12e27ae9fdfd390ede8643c4dd6b2656ba3e2256/e2e/65659143-8ae7-4026-91aa-8a2aa6f21603.md" "65659143-8ae7-4026-91aa-8a2aa6f21603.md"
Add-Link $ws1 "A4" "https://github.com/OpenLocalizationTest/oltest/blob/12e27ae9fdfd390ede8643c4dd6b2656ba3e2256/e2e/d4a6720e-54cf-49b4-a943-9d2147c24d82.md" "d4a6720e-54cf-49b4-a943-9d2147c24d82.md"
Add-Link $ws1 "A5" "https://github.com/OpenLocalizationTest/oltest/blob/7a2891d7c002b8b88399127f7118505501e6f1b8/e2e/5c4cc5a0-b7f6-4851-8ce5-f381df8f46f4.md" "5c4cc5a0-b7f6-4851-8ce5-f381df8f46f4.md"
Add-Link $ws1 "A6" "https://github.com/OpenLocalizationTest/oltest/blob/b2a7f48aeb67d425f95ce1a4edef38107ba415cf/e2e/ba776950-b98c-4ead-8535-7f56a0869957.md" "ba776950-b98c-4ead-8535-7f56a0869957.md"
Add-Link $ws1 "A7" "https://github.com/OpenLocalizationTest/oltest/blob/ff6286bfda74819fa80c209c072852480f52539c/e2e/c97306d9-c77a-49c6-abdf-eea21385d93f.md" "c97306d9-c77a-49c6-abdf-eea21385d93f.md"

# ---------------------------------------------------------------------
# Shared helper for the per-language detail sheets ("zh-cn" / "de-de")
# Columns: A Source File Name | B File Extension | C Status
#          D Latest Handoff File | E Latest Handoff Datetime
#          F Latest Target File | G Latest Handback File
#          H Latest Handback DateTime | I Handoff Reason
# ---------------------------------------------------------------------
function Update-DetailSheet($ws, $lang, $rows) {
    foreach ($r in $rows) {
        $row = $r.Row
        $ws.Range("A$row").Value2 = $r.A
        $ws.Range("B$row").Value2 = $r.B
        $ws.Range("C$row").Value2 = $r.C
        $ws.Range("D$row").Value2 = $r.D
        $ws.Range("E$row").NumberFormat = "yyyy-mm-dd HH:mm:ss"
        $ws.Range("E$row").Value2 = $r.E
        if ($r.ContainsKey("F")) { $ws.Range("F$row").Value2 = $r.F }
        if ($r.ContainsKey("G")) { $ws.Range("G$row").Value2 = $r.G }
        $ws.Range("H$row").Value2 = $r.H
        $ws.Range("I$row").Value2 = $r.I
    }

    # Rebuild hyperlinks (A, B, D, and F/G on row 2 only) from scratch.
    $ws.Hyperlinks.Delete() | Out-Null

    foreach ($r in $rows) {
        $row = $r.Row
        Add-Link $ws "A$row" $r.AUrl $r.A
        Add-Link $ws "B$row" $r.BUrl $r.B
        Add-Link $ws "D$row" $r.DUrl $r.D
        if ($r.ContainsKey("FUrl")) { Add-Link $ws "F$row" $r.FUrl $r.F }
        if ($r.ContainsKey("GUrl")) { Add-Link $ws "G$row" $r.GUrl $r.G }
    }
}

# ---------------------------------------------------------------------
# Sheet 2: "zh-cn"
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")

$zh = @(
    @{ Row=2
       A="8e57db03-c7a5-42e7-aa11-fa025c83e4e3.md"; AUrl="https://github.com/OpenLocalizationTest/oltest/blob/07144325e6973bba9c1363117d82640c93e928c3/e2e/8e57db03-c7a5-42e7-aa11-fa025c83e4e3.md"
       B=".md"; BUrl="https://github.com/OpenLocalizationTest/oltest/blob/07144325e6973bba9c1363117d82640c93e928c3/e2e/8e57db03-c7a5-42e7-aa11-fa025c83e4e3.md"
       C="Handed back: in sync with en-US"
       D="8e57db03-c7a5-42e7-aa11-fa025c83e4e3.953ebf2f90085ebd3591317ac0a873d6bb0ee08f.zh-cn.xlf"; DUrl="https://github.com/OpenLocalizationTestOrg/olhandoff/blob/44b07ec9398facf38268d5b6e7c84afdc4ccbe6f/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/8e57db03-c7a5-42e7-aa11-fa025c83e4e3.953ebf2f90085ebd3591317ac0a873d6bb0ee08f.zh-cn.xlf"
       E="2016-03-12 02:25:38"
       F="8e57db03-c7a5-42e7-aa11-fa025c83e4e3.md"; FUrl="https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/a155ce243e017cbc31cb5b56c069c8af7041438f/e2e/8e57db03-c7a5-42e7-aa11-fa025c83e4e3.md"
       G="8e57db03-c7a5-42e7-aa11-fa025c83e4e3.953ebf2f90085ebd3591317ac0a873d6bb0ee08f.zh-cn.xlf"; GUrl="https://github.com/OpenLocalizationTestOrg/olhandback/blob/d38f3dc6102df952c29516fb5aca52550c9c6cd4/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/8e57db03-c7a5-42e7-aa11-fa025c83e4e3.953ebf2f90085ebd3591317ac0a873d6bb0ee08f.zh-cn.xlf"
       H="2016-03-12 02:25:54"; I="Include" },

    @{ Row=3
       A="65659143-8ae7-4026-91aa-8a2aa6f21603.md"; AUrl="https://github.com/OpenLocalizationTest/oltest/blob/12e27ae9fdfd390ede8643c4dd6b2656ba3e2256/e2e/65659143-8ae7-4026-91aa-8a2aa6f21603.md"
       B=".md"; BUrl="https://github.com/OpenLocalizationTest/oltest/blob/12e27ae9fdfd390ede8643c4dd6b2656ba3e2256/e2e/65659143-8ae7-4026-91aa-8a2aa6f21603.md"
       C="In Translation"
       D="65659143-8ae7-4026-91aa-8a2aa6f21603.db72ae1645917d864d7c74ef47f05a7268fe1785.zh-cn.xlf"; DUrl="https://github.com/OpenLocalizationTestOrg/olhandoff/blob/80a7e720c911df923cc7825c03b4b272dde784bd/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/65659143-8ae7-4026-91aa-8a2aa6f21603.db72ae1645917d864d7c74ef47f05a7268fe1785.zh-cn.xlf"
       E="2016-03-12 02:26:23"
       H="0001-01-01 00:00:00"; I="Include" },

    @{ Row=4
       A="d4a6720e-54cf-49b4-a943-9d2147c24d82.md"; AUrl="https://github.com/OpenLocalizationTest/oltest/blob/12e27ae9fdfd390ede8643c4dd6b2656ba3e2256/e2e/d4a6720e-54cf-49b4-a943-9d2147c24d82.md"
       B=".md"; BUrl="https://github.com/OpenLocalizationTest/oltest/blob/12e27ae9fdfd390ede8643c4dd6b2656ba3e2256/e2e/d4a6720e-54cf-49b4-a943-9d2147c24d82.md"
       C="In Translation"
       D="d4a6720e-54cf-49b4-a943-9d2147c24d82.5719726605eb21ea9db609ebcbbd7de571612818.zh-cn.xlf"; DUrl="https://github.com/OpenLocalizationTestOrg/olhandoff/blob/80a7e720c911df923cc7825c03b4b272dde784bd/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/d4a6720e-54cf-49b4-a943-9d2147c24d82.5719726605eb21ea9db609ebcbbd7de571612818.zh-cn.xlf"
       E="2016-03-12 02:26:23"
       H="0001-01-01 00:00:00"; I="Include" },

    @{ Row=5
       A="5c4cc5a0-b7f6-4851-8ce5-f381df8f46f4.md"; AUrl="https://github.com/OpenLocalizationTest/oltest/blob/7a2891d7c002b8b88399127f7118505501e6f1b8/e2e/5c4cc5a0-b7f6-4851-8ce5-f381df8f46f4.md"
       B=".md"; BUrl="https://github.com/OpenLocalizationTest/oltest/blob/7a2891d7c002b8b88399127f7118505501e6f1b8/e2e/5c4cc5a0-b7f6-4851-8ce5-f381df8f46f4.md"
       C="Ready for handoff"
       D="5c4cc5a0-b7f6-4851-8ce5-f381df8f46f4.7a2891d7c002b8b88399127f7118505501e6f1b8.zh-cn.xlf"; DUrl="https://github.com/OpenLocalizationTestOrg/olhandoff/blob/7a2891d7c002b8b88399127f7118505501e6f1b8/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/5c4cc5a0-b7f6-4851-8ce5-f381df8f46f4.7a2891d7c002b8b88399127f7118505501e6f1b8.zh-cn.xlf"
       E="2016-03-12 02:26:56"
       H="0001-01-01 00:00:00"; I="Include" },

    @{ Row=6
       A="ba776950-b98c-4ead-8535-7f56a0869957.md"; AUrl="https://github.com/OpenLocalizationTest/oltest/blob/b2a7f48aeb67d425f95ce1a4edef38107ba415cf/e2e/ba776950-b98c-4ead-8535-7f56a0869957.md"
       B=".md"; BUrl="https://github.com/OpenLocalizationTest/oltest/blob/b2a7f48aeb67d425f95ce1a4edef38107ba415cf/e2e/ba776950-b98c-4ead-8535-7f56a0869957.md"
       C="Ready for handoff"
       D="ba776950-b98c-4ead-8535-7f56a0869957.9e4c20aa6e62045b5f15ef2371c067e03c88c415.zh-cn.xlf"; DUrl="https://github.com/OpenLocalizationTestOrg/olhandoff/blob/6073ef83578e17ab55c5ee22608642b2ace4ba40/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/ba776950-b98c-4ead-8535-7f56a0869957.9e4c20aa6e62045b5f15ef2371c067e03c88c415.zh-cn.xlf"
       E="2016-03-12 02:25:16"
       H="0001-01-01 00:00:00"; I="Include" },

    @{ Row=7
       A="c97306d9-c77a-49c6-abdf-eea21385d93f.md"; AUrl="https://github.com/OpenLocalizationTest/oltest/blob/ff6286bfda74819fa80c209c072852480f52539c/e2e/c97306d9-c77a-49c6-abdf-eea21385d93f.md"
       B=".md"; BUrl="https://github.com/OpenLocalizationTest/oltest/blob/ff6286bfda74819fa80c209c072852480f52539c/e2e/c97306d9-c77a-49c6-abdf-eea21385d93f.md"
       C="Ready for handoff"
       D="c97306d9-c77a-49c6-abdf-eea21385d93f.ff6286bfda74819fa80c209c072852480f52539c.zh-cn.xlf"; DUrl="https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ff6286bfda74819fa80c209c072852480f52539c/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/c97306d9-c77a-49c6-abdf-eea21385d93f.ff6286bfda74819fa80c209c072852480f52539c.zh-cn.xlf"
       E="2016-03-12 02:26:56"
       H="0001-01-01 00:00:00"; I="Include" }
)

Update-DetailSheet $ws2 "zh-cn" $zh

# ---------------------------------------------------------------------
# Sheet 3: "de-de"
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")

$de = @(
    @{ Row=2
       A="8e57db03-c7a5-42e7-aa11-fa025c83e4e3.md"; AUrl="https://github.com/OpenLocalizationTest/oltest/blob/07144325e6973bba9c1363117d82640c93e928c3/e2e/8e57db03-c7a5-42e7-aa11-fa025c83e4e3.md"
       B=".md"; BUrl="https://github.com/OpenLocalizationTest/oltest/blob/07144325e6973bba9c1363117d82640c93e928c3/e2e/8e57db03-c7a5-42e7-aa11-fa025c83e4e3.md"
       C="Handed back: in sync with en-US"
       D="8e57db03-c7a5-42e7-aa11-fa025c83e4e3.953ebf2f90085ebd3591317ac0a873d6bb0ee08f.de-de.xlf"; DUrl="https://github.com/OpenLocalizationTestOrg/olhandoff/blob/fc2c60b3104014aef3802feb71c06be339879ccf/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/8e57db03-c7a5-42e7-aa11-fa025c83e4e3.953ebf2f90085ebd3591317ac0a873d6bb0ee08f.de-de.xlf"
       E="2016-03-12 02:25:41"
       F="8e57db03-c7a5-42e7-aa11-fa025c83e4e3.md"; FUrl="https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/a2169c16969e3a613a23110295f963f56e2737fd/e2e/8e57db03-c7a5-42e7-aa11-fa025c83e4e3.md"
       G="8e57db03-c7a5-42e7-aa11-fa025c83e4e3.953ebf2f90085ebd3591317ac0a873d6bb0ee08f.de-de.xlf"; GUrl="https://github.com/OpenLocalizationTestOrg/olhandback/blob/8448e11c122b3ac3753e2ccdd2e360ea1418e434/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/8e57db03-c7a5-42e7-aa11-fa025c83e4e3.953ebf2f90085ebd3591317ac0a873d6bb0ee08f.de-de.xlf"
       H="2016-03-12 02:26:00"; I="Include" },

    @{ Row=3
       A="65659143-8ae7-4026-91aa-8a2aa6f21603.md"; AUrl="https://github.com/OpenLocalizationTest/oltest/blob/12e27ae9fdfd390ede8643c4dd6b2656ba3e2256/e2e/65659143-8ae7-4026-91aa-8a2aa6f21603.md"
       B=".md"; BUrl="https://github.com/OpenLocalizationTest/oltest/blob/12e27ae9fdfd390ede8643c4dd6b2656ba3e2256/e2e/65659143-8ae7-4026-91aa-8a2aa6f21603.md"
       C="In Translation"
       D="65659143-8ae7-4026-91aa-8a2aa6f21603.db72ae1645917d864d7c74ef47f05a7268fe1785.de-de.xlf"; DUrl="https://github.com/OpenLocalizationTestOrg/olhandoff/blob/cb65ff2d196923071b41e29578ea00feb3060581/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/65659143-8ae7-4026-91aa-8a2aa6f21603.db72ae1645917d864d7c74ef47f05a7268fe1785.de-de.xlf"
       E="2016-03-12 02:26:26"
       H="0001-01-01 00:00:00"; I="Include" },

    @{ Row=4
       A="d4a6720e-54cf-49b4-a943-9d2147c24d82.md"; AUrl="https://github.com/OpenLocalizationTest/oltest/blob/12e27ae9fdfd390ede8643c4dd6b2656ba3e2256/e2e/d4a6720e-54cf-49b4-a943-9d2147c24d82.md"
       B=".md"; BUrl="https://github.com/OpenLocalizationTest/oltest/blob/12e27ae9fdfd390ede8643c4dd6b2656ba3e2256/e2e/d4a6720e-54cf-49b4-a943-9d2147c24d82.md"
       C="In Translation"
       D="d4a6720e-54cf-49b4-a943-9d2147c24d82.5719726605eb21ea9db609ebcbbd7de571612818.de-de.xlf"; DUrl="https://github.com/OpenLocalizationTestOrg/olhandoff/blob/cb65ff2d196923071b41e29578ea00feb3060581/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/d4a6720e-54cf-49b4-a943-9d2147c24d82.5719726605eb21ea9db609ebcbbd7de571612818.de-de.xlf"
       E="2016-03-12 02:26:26"
       H="0001-01-01 00:00:00"; I="Include" },

    @{ Row=5
       A="5c4cc5a0-b7f6-4851-8ce5-f381df8f46f4.md"; AUrl="https://github.com/OpenLocalizationTest/oltest/blob/7a2891d7c002b8b88399127f7118505501e6f1b8/e2e/5c4cc5a0-b7f6-4851-8ce5-f381df8f46f4.md"
       B=".md"; BUrl="https://github.com/OpenLocalizationTest/oltest/blob/7a2891d7c002b8b88399127f7118505501e6f1b8/e2e/5c4cc5a0-b7f6-4851-8ce5-f381df8f46f4.md"
       C="Ready for handoff"
       D="5c4cc5a0-b7f6-4851-8ce5-f381df8f46f4.7a2891d7c002b8b88399127f7118505501e6f1b8.de-de.xlf"; DUrl="https://github.com/OpenLocalizationTestOrg/olhandoff/blob/7a2891d7c002b8b88399127f7118505501e6f1b8/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/5c4cc5a0-b7f6-4851-8ce5-f381df8f46f4.7a2891d7c002b8b88399127f7118505501e6f1b8.de-de.xlf"
       E="2016-03-12 02:26:59"
       H="0001-01-01 00:00:00"; I="Include" },

    @{ Row=6
       A="ba776950-b98c-4ead-8535-7f56a0869957.md"; AUrl="https://github.com/OpenLocalizationTest/oltest/blob/b2a7f48aeb67d425f95ce1a4edef38107ba415cf/e2e/ba776950-b98c-4ead-8535-7f56a0869957.md"
       B=".md"; BUrl="https://github.com/OpenLocalizationTest/oltest/blob/b2a7f48aeb67d425f95ce1a4edef38107ba415cf/e2e/ba776950-b98c-4ead-8535-7f56a0869957.md"
       C="Ready for handoff"
       D="ba776950-b98c-4ead-8535-7f56a0869957.9e4c20aa6e62045b5f15ef2371c067e03c88c415.de-de.xlf"; DUrl="https://github.com/OpenLocalizationTestOrg/olhandoff/blob/2d3103129974ae24ddfd840a0bb64e093e177076/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/ba776950-b98c-4ead-8535-7f56a0869957.9e4c20aa6e62045b5f15ef2371c067e03c88c415.de-de.xlf"
       E="2016-03-12 02:25:19"
       H="0001-01-01 00:00:00"; I="Include" },

    @{ Row=7
       A="c97306d9-c77a-49c6-abdf-eea21385d93f.md"; AUrl="https://github.com/OpenLocalizationTest/oltest/blob/ff6286bfda74819fa80c209c072852480f52539c/e2e/c97306d9-c77a-49c6-abdf-eea21385d93f.md"
       B=".md"; BUrl="https://github.com/OpenLocalizationTest/oltest/blob/ff6286bfda74819fa80c209c072852480f52539c/e2e/c97306d9-c77a-49c6-abdf-eea21385d93f.md"
       C="Ready for handoff"
       D="c97306d9-c77a-49c6-abdf-eea21385d93f.ff6286bfda74819fa80c209c072852480f52539c.de-de.xlf"; DUrl="https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ff6286bfda74819fa80c209c072852480f52539c/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/c97306d9-c77a-49c6-abdf-eea21385d93f.ff6286bfda74819fa80c209c072852480f52539c.de-de.xlf"
       E="2016-03-12 02:26:59"
       H="0001-01-01 00:00:00"; I="Include" }
)

Update-DetailSheet $ws3 "de-de" $de

$wb.Save()
